# Update countries & provincias Spain
#
# The data refresh bumps the "last updated" timestamp and a number of
# per-country case counts. Because the sheet is kept sorted in descending
# order by "Casos totales" (column B), two countries whose totals grew
# climbed past their neighbours:
#   - Colombia overtook Japon / Austria / Kuwait (rows 40-43)
#   - Guyana overtook Bermudas / Camboya (rows 162-164)
# so those rows are rewritten in their new order with refreshed figures,
# while a handful of other rows just get refreshed totals in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header / timestamp (row 1) -----------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 19 de Mayo de 2020 a las 00:05"

# --- Plain data refreshes (ranking unchanged) ----------------------------
# Estados Unidos (row 4)
$ws.Cells.Item(4,2).Value = 1544874
$ws.Cells.Item(4,3).Value = 17210
$ws.Cells.Item(4,4).Value = 353172
$ws.Cells.Item(4,5).Value = 1099971
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 753
$ws.Cells.Item(4,8).Value = 91731

# Alemania (row 11)
$ws.Cells.Item(11,2).Value = 177289
$ws.Cells.Item(11,3).Value = 638
$ws.Cells.Item(11,4).Value = 154600
$ws.Cells.Item(11,5).Value = 14566
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 74
$ws.Cells.Item(11,8).Value = 8123

# Reunion (row 131) - only B/C/D/E/F change
$ws.Cells.Item(131,2).Value = 446
$ws.Cells.Item(131,3).Value = 3
$ws.Cells.Item(131,4).Value = 354
$ws.Cells.Item(131,5).Value = 92
$ws.Cells.Item(131,6).Value = 0

# Santo Tome y Principe (row 149) - only B/C/D/E/F change
$ws.Cells.Item(149,2).Value = 246
$ws.Cells.Item(149,3).Value = 11
$ws.Cells.Item(149,4).Value = 4
$ws.Cells.Item(149,5).Value = 235
$ws.Cells.Item(149,6).Value = 0

# --- Colombia climbs above Japon / Austria / Kuwait (rows 40-43) --------
$ws.Cells.Item(40,1).Value = "Colombia"
$ws.Cells.Item(40,2).Value = 16295
$ws.Cells.Item(40,3).Value = 721
$ws.Cells.Item(40,4).Value = 3903
$ws.Cells.Item(40,5).Value = 11800
$ws.Cells.Item(40,6).Value = 0
$ws.Cells.Item(40,7).Value = 18
$ws.Cells.Item(40,8).Value = 592

$ws.Cells.Item(41,1).Value = "Japon"
$ws.Cells.Item(41,2).Value = 16285
$ws.Cells.Item(41,3).Value = 0
$ws.Cells.Item(41,4).Value = 11153
$ws.Cells.Item(41,5).Value = 4388
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 0
$ws.Cells.Item(41,8).Value = 744

$ws.Cells.Item(42,1).Value = "Austria"
$ws.Cells.Item(42,2).Value = 16269
$ws.Cells.Item(42,3).Value = 27
$ws.Cells.Item(42,4).Value = 14614
$ws.Cells.Item(42,5).Value = 1026
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 0
$ws.Cells.Item(42,8).Value = 629

$ws.Cells.Item(43,1).Value = "Kuwait"
$ws.Cells.Item(43,2).Value = 15691
$ws.Cells.Item(43,3).Value = 841
$ws.Cells.Item(43,4).Value = 4339
$ws.Cells.Item(43,5).Value = 11234
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 6
$ws.Cells.Item(43,8).Value = 118

# --- Guyana climbs above Bermudas / Camboya (rows 162-164) --------------
$ws.Cells.Item(162,1).Value = "Guyana"
$ws.Cells.Item(162,2).Value = 124
$ws.Cells.Item(162,3).Value = 7
$ws.Cells.Item(162,4).Value = 44
$ws.Cells.Item(162,5).Value = 70
$ws.Cells.Item(162,6).Value = 0
$ws.Cells.Item(162,7).Value = 0
$ws.Cells.Item(162,8).Value = 10

$ws.Cells.Item(163,1).Value = "Bermudas"
$ws.Cells.Item(163,2).Value = 123
$ws.Cells.Item(163,3).Value = 0
$ws.Cells.Item(163,4).Value = 73
$ws.Cells.Item(163,5).Value = 41
$ws.Cells.Item(163,6).Value = 0
$ws.Cells.Item(163,7).Value = 0
$ws.Cells.Item(163,8).Value = 9

$ws.Cells.Item(164,1).Value = "Camboya"
$ws.Cells.Item(164,2).Value = 122
$ws.Cells.Item(164,3).Value = 0
$ws.Cells.Item(164,4).Value = 122
$ws.Cells.Item(164,5).Value = 0
$ws.Cells.Item(164,6).Value = 0
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 0
